# Applies the "LinuxForHealth" re-branding update to the
# StructureDefinition-source-record-type workbook:
#   - Metadata sheet: URL, Version, Date and Publisher values are updated
#   - Elements sheet: the (now stale) Constraint(s) note on the root
#     "Extension" row is cleared out

$wb = $excel.ActiveWorkbook

# --- Metadata sheet -------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/source-record-type"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- Elements sheet ---------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

# Clear the Constraint(s) cell for the root "Extension" row (row 2)
$elements.Range("AI2").ClearContents()

# The Extension.url row carries the same canonical URL as its Fixed Value;
# keep it in sync with the updated URL above
$elements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/source-record-type"
